$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value  = "Done"
$ws.Range("C4").Value  = "Done"
$ws.Range("C5").Value  = "Half Done"
$ws.Range("C6").Value  = "Done"
$ws.Range("C7").Value  = "Done"
$ws.Range("C8").Value  = "Almost Done"
$ws.Range("C9").Value  = "Done"
$ws.Range("C10").Value = "Pending"
$ws.Range("C11").Value = "Hold - Admin Page"
